# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the existing columns (bold + border, style index 1) and
# filling the data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it keeps the same bold/border style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
